$d = $word.ActiveDocument

$replacements = @(
    @("283÷9=", "326÷4="),
    @("770÷8=", "356÷4="),
    @("693÷3=", "987÷4="),
    @("338÷6=", "141÷2="),
    @("589÷9=", "577÷8="),
    @("699÷6=", "752÷5="),
    @("793÷4=", "130÷6="),
    @("494÷9=", "213÷5="),
    @("452÷6=", "248÷5="),
    @("709÷9=", "797÷4="),
    @("101÷8=", "556÷5="),
    @("420÷3=", "478÷8="),
    @("659÷7=", "541÷5="),
    @("733÷5=", "933÷8="),
    @("661÷5=", "180÷2="),
    @("466÷2=", "380÷6="),
    @("846÷8=", "249÷2="),
    @("278÷3=", "112÷2="),
    @("670÷2=", "641÷7="),
    @("735÷7=", "819÷4="),
    @("131÷5=", "229÷9="),
    @("610÷5=", "695÷2="),
    @("339÷3=", "782÷2="),
    @("336÷2=", "154÷7="),
    @("843÷2=", "693÷4=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
